$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet2: "M"/"F" -> "Male"/"Female", add two new rows for "person3",
# and wire up the matching mailto hyperlinks.
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Sheet2")

# Existing gender abbreviations become spelled-out values.
$ws2.Range("A2").Value = "Male"
$ws2.Range("A3").Value = "Female"
$ws2.Range("A4").Value = "Male"

# New row 5: person3 (female)
$ws2.Range("A5").Value = "Female"
$ws2.Range("B5").Value = "person3"
$ws2.Range("D5").Value = "person3"
$ws2.Range("E5").Value = "person3"
$ws2.Range("F5").Value = "person"
$ws2.Range("G5").Value = "person"
$ws2.Range("H5").Value = "person"
$ws2.Range("I5").Value = "somewhere"
$ws2.Hyperlinks.Add($ws2.Range("C5"), "mailto:person3@email.com", [Type]::Missing, [Type]::Missing, "person3@email.com")

# New row 6: duplicate of row 5
$ws2.Range("A6").Value = "Female"
$ws2.Range("B6").Value = "person3"
$ws2.Range("D6").Value = "person3"
$ws2.Range("E6").Value = "person3"
$ws2.Range("F6").Value = "person"
$ws2.Range("G6").Value = "person"
$ws2.Range("H6").Value = "person"
$ws2.Range("I6").Value = "somewhere"
$ws2.Hyperlinks.Add($ws2.Range("C6"), "mailto:person3@email.com", [Type]::Missing, [Type]::Missing, "person3@email.com")

# ---------------------------------------------------------------------
# Sheet4: stray "yespassword" entries were lost/retyped as "mypassword"
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Sheet4")
$ws4.Range("C3").Value = "mypassword"
$ws4.Range("A4").Value = "mypassword"

# ---------------------------------------------------------------------
# Selection / active-sheet bookkeeping
# ---------------------------------------------------------------------
$ws4.Range("A4").Select()

$ws2.Activate()
$ws2.Range("B6").Select()
